$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.533.54'
$ws.Range('E2').Value = '  -2.25%  '
$ws.Range('D3').Value = '3.681.64'
$ws.Range('E3').Value = '  -3.05%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '614.28'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.08'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.70%  '
$ws.Range('D7').Value = '3.678.57'
$ws.Range('E7').Value = '  -3.26%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -3.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.25'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.29%  '
$ws.Range('E12').Value = '  -5.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '39.84'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.78%  '
$ws.Range('E14').Value = '  -4.32%  '
$ws.Range('D15').Value = '4.299.66'
$ws.Range('E15').Value = '  -3.24%  '
$ws.Range('D16').Value = '3.683.00'
$ws.Range('E16').Value = '  -3.26%  '
$ws.Range('D17').Value = '69.506.10'
$ws.Range('E17').Value = '  -2.44%  '
$ws.Range('E18').Value = '  -1.88%  '
$ws.Range('E19').Value = '  -1.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.33'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '499.05'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.96%  '
$ws.Range('E22').Value = '  -3.93%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.714'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.16%  '
$ws.Range('E24').Value = '  -1.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.16'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.26'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.91'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.52%  '
$ws.Range('E28').Value = '  +0.97%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.43'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.72%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.88'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.94'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '30.01'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -8.03%  '
$ws.Range('E34').Value = '  -2.46%  '
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('E36').Value = '  -1.82%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.02'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.137'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.338'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.96%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '49.85'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.77%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.05'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -8.29%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.92'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.56%  '
$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '428.29'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.96'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.55'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.94%  '
$ws.Range('D46').Value = '2.926.55'
$ws.Range('E46').Value = '  -7.52%  '
$ws.Range('E47').Value = '  -3.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '27.27'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '136.24'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.43'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.55%  '
